$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, avoiding Excels auto number
# coercion for numeric-looking strings (e.g. "224.60" -> 224.6), while
# keeping the cells original (unstyled) formatting by clearing the
# temporary text NumberFormat we apply to force the literal string.
function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "95.696.14"
Set-TextValue "E2" "  -1.90%  "
Set-TextValue "D3" "3.630.34"
Set-TextValue "E3" "  -2.39%  "
Set-TextValue "D4" "2.77"
Set-TextValue "E4" "  +27.12%  "
Set-TextValue "E5" "  +0.11%  "
Set-TextValue "D6" "224.60"
Set-TextValue "E6" "  -5.70%  "
Set-TextValue "D7" "644.57"
Set-TextValue "E7" "  -1.78%  "
Set-TextValue "D8" "0.423"
Set-TextValue "E8" "  -4.72%  "
Set-TextValue "D9" "1.20"
Set-TextValue "E9" "  +4.92%  "
Set-TextValue "E10" "  +0.01%  "
Set-TextValue "D11" "3.627.65"
Set-TextValue "E11" "  -2.35%  "
Set-TextValue "D12" "51.46"
Set-TextValue "E12" "  +12.97%  "
Set-TextValue "E13" "  +5.73%  "
Set-TextValue "D14" "0.0000295"
Set-TextValue "E14" "  -5.37%  "
Set-TextValue "D15" "6.54"
Set-TextValue "E15" "  -4.35%  "
Set-TextValue "D16" "4.307.14"
Set-TextValue "E16" "  -2.45%  "
Set-TextValue "D17" "24.97"
Set-TextValue "E17" "  +31.86%  "
Set-TextValue "D18" "95.482.54"
Set-TextValue "E18" "  -1.80%  "
Set-TextValue "D19" "9.33"
Set-TextValue "E19" "  +5.16%  "
Set-TextValue "D20" "13.93"
Set-TextValue "E20" "  +6.37%  "
Set-TextValue "D21" "3.623.51"
Set-TextValue "E21" "  -2.54%  "
Set-TextValue "D22" "0.314"
Set-TextValue "E22" "  +46.08%  "
Set-TextValue "D23" "0.540"
Set-TextValue "E23" "  -0.06%  "
Set-TextValue "D24" "136.26"
Set-TextValue "E24" "  +14.56%  "
Set-TextValue "D25" "534.12"
Set-TextValue "E25" "  +0.57%  "
Set-TextValue "D26" "3.33"
Set-TextValue "E26" "  -4.30%  "
Set-TextValue "D27" "7.13"
Set-TextValue "E27" "  +3.12%  "
Set-TextValue "D28" "0.0000203"
Set-TextValue "E28" "  -9.05%  "
Set-TextValue "D29" "13.56"
Set-TextValue "E29" "  +0.85%  "
Set-TextValue "D30" "3.798.09"
Set-TextValue "E30" "  -3.09%  "
Set-TextValue "D31" "13.67"
Set-TextValue "E31" "  +6.62%  "
Set-TextValue "D32" "3.20"
Set-TextValue "E32" "  +5.11%  "
Set-TextValue "E33" "  +0.04%  "
Set-TextValue "D34" "0.649"
Set-TextValue "E34" "  +7.02%  "
Set-TextValue "E35" "  +2.74%  "
Set-TextValue "D36" "33.87"
Set-TextValue "E36" "  +2.32%  "
Set-TextValue "E37" "  -4.94%  "
Set-TextValue "D38" "0.999"
Set-TextValue "E38" "  +0.75%  "
Set-TextValue "D39" "0.0565"
Set-TextValue "E39" "  +22.18%  "
Set-TextValue "B40" "RenderToken"
Set-TextValue "C40" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D40" "8.57"
Set-TextValue "E40" "  -1.83%  "
Set-TextValue "B41" "USDe"
Set-TextValue "C41" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  -0.04%  "
Set-TextValue "D42" "601.36"
Set-TextValue "E42" "  -6.14%  "
Set-TextValue "D43" "7.24"
Set-TextValue "E43" "  +4.94%  "
Set-TextValue "D44" "0.507"
Set-TextValue "E45" "  +5.53%  "
Set-TextValue "D46" "41.36"
Set-TextValue "E46" "  -0.05%  "
Set-TextValue "D47" "2.03"
Set-TextValue "E47" "  +0.30%  "
Set-TextValue "E48" "  -7.14%  "
Set-TextValue "D49" "9.40"
Set-TextValue "E49" "  +4.77%  "
Set-TextValue "D50" "238.13"
Set-TextValue "E50" "  +13.68%  "
Set-TextValue "D51" "2.38"
Set-TextValue "E51" "  -1.04%  "
